$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the phi (φ) metric values to 0
$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 0
$ws.Range("B4").Value = 0
$ws.Range("B5").Value = 0

# Update the (ms) timing metric values with new measurements
$ws.Range("B7").Value = 0.05704092979431152
$ws.Range("B8").Value = 0.2069904804229736
$ws.Range("B9").Value = 0.1339104175567627
$ws.Range("B10").Value = 0.2220222949981689
